$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.446.72'
$ws.Range('E2').Value = '  -1.24%  '
$ws.Range('D3').Value = '1.685.77'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('D4').Value = '''0.9995'
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = '''316.58'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').Value = '''0.9997'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').Value = '''0.3882'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('D8').Value = '''0.4013'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').Value = '''1.484'
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').Value = '''0.9999'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').Value = '''52.45'
$ws.Range('E11').Value = '  -2.89%  '
$ws.Range('D12').Value = '''0.08757'
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '''26.15'
$ws.Range('E13').Value = '  +11.99%  '
$ws.Range('D14').Value = '''7.488'
$ws.Range('E14').Value = '  +3.75%  '
$ws.Range('D15').Value = '''8.008'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '''0.00001344'
$ws.Range('E16').Value = '  +1.25%  '
$ws.Range('D17').Value = '1.664.10'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').Value = '''97.73'
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('D19').Value = '''0.07217'
$ws.Range('E19').Value = '  +3.04%  '
$ws.Range('D20').Value = '''19.72'
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').Value = '''7.256'
$ws.Range('E21').Value = '  +3.52%  '
$ws.Range('D22').Value = '''0.9998'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '''14.17'
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('D24').Value = '24.428.29'
$ws.Range('E24').Value = '  -1.32%  '
$ws.Range('D25').Value = '''3.024'
$ws.Range('E25').Value = '  -6.83%  '
$ws.Range('D26').Value = '''2.346'
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('D27').Value = '''22.57'
$ws.Range('E27').Value = '  -0.67%  '
$ws.Range('D28').Value = '''168.08'
$ws.Range('E28').Value = '  +4.39%  '
$ws.Range('D29').Value = '''8.591'
$ws.Range('E29').Value = '  +10.96%  '
$ws.Range('D30').Value = '''5.354'
$ws.Range('E30').Value = '  +3.67%  '
$ws.Range('D31').Value = '''138.18'
$ws.Range('E31').Value = '  +1.29%  '
$ws.Range('D32').Value = '1.849.54'
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('D33').Value = '''0.08766'
$ws.Range('E33').Value = '  +0.36%  '
$ws.Range('D34').Value = '''7.325'
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('D35').Value = '''1.052'
$ws.Range('E35').Value = '  -1.77%  '
$ws.Range('D36').Value = '''0.03004'
$ws.Range('E36').Value = '  +9.98%  '
$ws.Range('D37').Value = '''1.976'
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').Value = '''0.2750'
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('D39').Value = '''10.80'
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('D40').Value = '''0.09133'
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').Value = '''0.8023'
$ws.Range('E41').Value = '  +4.74%  '
$ws.Range('D42').Value = '''14.07'
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').Value = '''1.474'
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('D44').Value = '''17.70'
$ws.Range('E44').Value = '  +12.17%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '''2.610'
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.7220'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').Value = '''4.265'
$ws.Range('E47').Value = '  +1.24%  '
$ws.Range('D48').Value = '''1.400'
$ws.Range('E48').Value = '  +7.01%  '
$ws.Range('D49').Value = '''0.9993'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').Value = '''139.22'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('D51').Value = '''0.08050'
$ws.Range('E51').Value = '  +0.90%  '

Write-Output "Updated cryptos list"
